$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 1135442.8
$ws.Range("I70").Value = 2551771.8
$ws.Range("J70").Value = 2379.6
$ws.Range("K70").Value = 7655315.399999999
$ws.Range("L70").Value = 7138.799999999999
$ws.Range("M70").Value = -7655045.399999999
$ws.Range("N70").Value = -7678.799999999999
# Row 73
$ws.Range("H73").Value = 1135442.8
$ws.Range("I73").Value = 2551771.8
$ws.Range("J73").Value = 2379.6
$ws.Range("K73").Value = 7655315.399999999
$ws.Range("L73").Value = 7138.799999999999
$ws.Range("M73").Value = -7654379.399999999
$ws.Range("N73").Value = -9010.799999999999
# Row 86
$ws.Range("H86").Value = 11827780
$ws.Range("I86").Value = 4642.25
$ws.Range("K86").Value = 4642.25
$ws.Range("M86").Value = -3519.25
# Row 88
$ws.Range("H88").Value = 2565.3809
$ws.Range("J88").Value = 2653.2666
$ws.Range("L88").Value = 2653.2666
$ws.Range("N88").Value = -3465.2666
# Row 89
$ws.Range("H89").Value = 11827780
$ws.Range("I89").Value = 4642.25
$ws.Range("K89").Value = 23211.25
$ws.Range("M89").Value = -17595.25
# Row 91
$ws.Range("H91").Value = 2565.3809
$ws.Range("J91").Value = 2653.2666
$ws.Range("L91").Value = 2653.2666
$ws.Range("N91").Value = -5461.2666
# Row 132
$ws.Range("H132").Value = 4728.5835
$ws.Range("I132").Value = 4148.8125
$ws.Range("K132").Value = 12446.4375
$ws.Range("M132").Value = -9916.4375
# Row 137
$ws.Range("H137").Value = 24393072
# Row 138
$ws.Range("H138").Value = 2256.2163
$ws.Range("I138").Value = 1647.3478
$ws.Range("J138").Value = 3256.5
$ws.Range("K138").Value = 4942.0434
$ws.Range("L138").Value = 9769.5
$ws.Range("M138").Value = 197.9565999999995
$ws.Range("N138").Value = -20049.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = $null
# Row 30
$ws.Range("H30").Value = 1887.25
$ws.Range("I30").Value = 2183
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 2183
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = -2033
$ws.Range("N30").Value = -1300
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null
# Row 132
$ws.Range("H132").Value = 2013.8948
$ws.Range("I132").Value = 2017.8667
$ws.Range("K132").Value = 6053.6001
$ws.Range("M132").Value = -3523.6001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 14032
$ws.Range("I20").Value = 7419.154
$ws.Range("K20").Value = 7419.154
$ws.Range("M20").Value = -7172.154
# Row 22
$ws.Range("H22").Value = 174.5
$ws.Range("I22").Value = 199
$ws.Range("K22").Value = 199
$ws.Range("M22").Value = -26
# Row 36
$ws.Range("H36").Value = 2566
$ws.Range("I36").Value = 2566
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2566
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2032
$ws.Range("N36").Value = $null
# Row 86
$ws.Range("H86").Value = 41669900
$ws.Range("I86").Value = 50002500
$ws.Range("J86").Value = 6902
$ws.Range("K86").Value = 50002500
$ws.Range("L86").Value = 6902
$ws.Range("M86").Value = -50001377
$ws.Range("N86").Value = -9148
# Row 89
$ws.Range("H89").Value = 41669900
$ws.Range("I89").Value = 50002500
$ws.Range("J89").Value = 6902
$ws.Range("K89").Value = 250012500
$ws.Range("L89").Value = 34510
$ws.Range("M89").Value = -250006884
$ws.Range("N89").Value = -45742
# Row 126
$ws.Range("H126").Value = 64171
$ws.Range("J126").Value = 64171
$ws.Range("L126").Value = 64171
$ws.Range("N126").Value = -74051
# Row 134
$ws.Range("H134").Value = 1014.4
$ws.Range("I134").Value = 1014.4
$ws.Range("K134").Value = 3043.2
$ws.Range("M134").Value = -508.1999999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 4900
$ws.Range("I62").Value = 4900
$ws.Range("K62").Value = 4900
$ws.Range("M62").Value = -4276
# Row 65
$ws.Range("H65").Value = 4900
$ws.Range("I65").Value = 4900
$ws.Range("K65").Value = 24500
$ws.Range("M65").Value = -21380
# Row 132
$ws.Range("H132").Value = 1591.5483
$ws.Range("I132").Value = 1492.0769
$ws.Range("K132").Value = 4476.2307
$ws.Range("M132").Value = -1946.2307
# Row 134
$ws.Range("H134").Value = 2358.111
$ws.Range("J134").Value = 3711
$ws.Range("L134").Value = 11133
$ws.Range("N134").Value = -16203
# Row 141
$ws.Range("H141").Value = 105706.25
$ws.Range("J141").Value = 105706.25
$ws.Range("L141").Value = 105706.25
$ws.Range("N141").Value = -116066.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 280.2
$ws.Range("I36").Value = 100.25
$ws.Range("J36").Value = 1000
$ws.Range("K36").Value = 300.75
$ws.Range("L36").Value = 3000
$ws.Range("M36").Value = -131.75
$ws.Range("N36").Value = -3338

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 6244.926
$ws.Range("I132").Value = 5405.857
$ws.Range("J132").Value = 7148.5386
$ws.Range("K132").Value = 16217.571
$ws.Range("L132").Value = 21445.6158
$ws.Range("M132").Value = -13687.571
$ws.Range("N132").Value = -26505.6158

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2860.75
$ws.Range("I40").Value = 2616.3
$ws.Range("K40").Value = 2616.3
$ws.Range("M40").Value = -2480.3
# Row 68
$ws.Range("H68").Value = 3509.6
$ws.Range("I68").Value = 4212.1665
$ws.Range("J68").Value = 2455.75
$ws.Range("K68").Value = 4212.1665
$ws.Range("L68").Value = 2455.75
$ws.Range("M68").Value = -3463.1665
$ws.Range("N68").Value = -3953.75
# Row 71
$ws.Range("H71").Value = 3509.6
$ws.Range("I71").Value = 4212.1665
$ws.Range("J71").Value = 2455.75
$ws.Range("K71").Value = 21060.8325
$ws.Range("L71").Value = 12278.75
$ws.Range("M71").Value = -17316.8325
$ws.Range("N71").Value = -19766.75
# Row 82
$ws.Range("H82").Value = 3969.5
$ws.Range("I82").Value = 2734.4167
$ws.Range("K82").Value = 2734.4167
$ws.Range("M82").Value = -2373.4167
# Row 85
$ws.Range("H85").Value = 3969.5
$ws.Range("I85").Value = 2734.4167
$ws.Range("K85").Value = 2734.4167
$ws.Range("M85").Value = -1486.4167

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 17
$ws.Range("H17").Value = 500502
$ws.Range("I17").Value = 500502
$ws.Range("K17").Value = 500502
$ws.Range("M17").Value = -500330
# Row 61
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1708
$ws.Range("N61").Value = $null
# Row 62
$ws.Range("H62").Value = 93092230
$ws.Range("I62").Value = 3974587.2
$ws.Range("J62").Value = 200033390
$ws.Range("K62").Value = 3974587.2
$ws.Range("L62").Value = 200033390
$ws.Range("M62").Value = -3973963.2
$ws.Range("N62").Value = -200034638
# Row 65
$ws.Range("H65").Value = 93092230
$ws.Range("I65").Value = 3974587.2
$ws.Range("J65").Value = 200033390
$ws.Range("K65").Value = 19872936
$ws.Range("L65").Value = 1000166950
$ws.Range("M65").Value = -19869816
$ws.Range("N65").Value = -1000173190
# Row 81
$ws.Range("H81").Value = 9533474
$ws.Range("I81").Value = 7848.4
$ws.Range("J81").Value = 18193132
$ws.Range("K81").Value = 15696.8
$ws.Range("L81").Value = 36386264
$ws.Range("M81").Value = -14635.8
$ws.Range("N81").Value = -36388386
# Row 84
$ws.Range("H84").Value = 9533474
$ws.Range("I84").Value = 7848.4
$ws.Range("J84").Value = 18193132
$ws.Range("K84").Value = 78484
$ws.Range("L84").Value = 181931320
$ws.Range("M84").Value = -73180
$ws.Range("N84").Value = -181941928
# Row 100
$ws.Range("H100").Value = 3620.75
$ws.Range("I100").Value = 4051.4707
$ws.Range("K100").Value = 8102.9414
$ws.Range("M100").Value = -7561.9414
# Row 132
$ws.Range("H132").Value = 4807
$ws.Range("I132").Value = 4309.1875
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 12927.5625
$ws.Range("L132").Value = 19200
$ws.Range("M132").Value = -10397.5625
$ws.Range("N132").Value = -24260
# Row 136
$ws.Range("H136").Value = 5607
$ws.Range("I136").Value = 2917.3333
$ws.Range("J136").Value = 7624.25
$ws.Range("K136").Value = 8751.999899999999
$ws.Range("L136").Value = 22872.75
$ws.Range("M136").Value = -6201.999899999999
$ws.Range("N136").Value = -27972.75
